# sefira.xlsx — "fixed nav color, added all days"
#
# For every day-of-the-omer row (1..48) this adds:
#   E<n> = ="{hebrew:'"&A<n>&"',"          (shared formula group si=0, rows 2:48)
#   F<n> = ="english:'"&B<n>&"'},"         (shared formula group si=1, rows 2:48)
#   G<n> = the literal text "{hebrew:'<hebrew text>',english:'<english text>'},"
#          (a plain string value — becomes a new shared-string entry, exactly
#          the concatenation of what E<n> and F<n> compute)
#
# then selects G1:G48, matching the saved workbook's cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E: ="{hebrew:'"&A<n>&"'," ------------------------------------
$ws.Range("E1").Formula = "=""{hebrew:'""&A1&""',"""
$ws.Range("E2:E48").Formula = "=""{hebrew:'""&A2&""',"""

# --- Column F: ="english:'"&B<n>&"'}," -----------------------------------
$ws.Range("F1").Formula = "=""english:'""&B1&""'},"""
$ws.Range("F2:F48").Formula = "=""english:'""&B2&""'},"""

# --- Column G: literal combined string, one new shared string per row ---
for ($r = 1; $r -le 48; $r++) {
    $hebrew = $ws.Cells.Item($r, 1).Value2
    $english = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 7).Value = "{hebrew:'" + $hebrew + "',english:'" + $english + "'},"
}

# Match the saved selection.
$null = $ws.Range("G1:G48").Select()
